# "Dados de seguranca analisados" - four extra SAEB/IDEB analysis columns
# (F:nota_saeb_matematica, G:nota_saeb_lingua_portuguesa,
#  H:nota_saeb_media_padronizada, I:ideb) had their column widths widened
# to fit the analyzed data, and the view was scrolled/re-selected further
# down the sheet (row 9 selected, view scrolled toward row 112).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen the newly analyzed data columns (F, G, H, I) to fit their content.
# ColumnWidth is expressed in "characters" and Excel quantizes it to the
# nearest pixel internally, so the inputs below are chosen to land as close
# as possible to the authored widths of 21.82 / 28.99 / 25.57 / 32.41.
$ws.Columns.Item(6).ColumnWidth = 21.0
$ws.Columns.Item(7).ColumnWidth = 28.166666666666668
$ws.Columns.Item(8).ColumnWidth = 24.666666666666668
$ws.Columns.Item(9).ColumnWidth = 31.5

# Move the selection to row 9 while reviewing the newly analyzed columns...
[void]$ws.Rows(9).Select()

# ...and scroll the window further down the data (toward row 112).
$excel.ActiveWindow.ScrollRow = 112
$excel.ActiveWindow.ScrollColumn = 1
